$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume/hour data for Thu Feb 16 19:05:57 UTC 2023 run.
# Each touched cell is first forced to text format ("@") so that
# numeric-looking values (prices, percentages, hour strings) are stored
# exactly as text, avoiding floating point re-interpretation and
# preserving formatting such as leading/trailing zeros and the "%" sign.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "320.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "5.88%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "19"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "49.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "11.83%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "19"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.308"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.02%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "19"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08065"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.20%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "19"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.595"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.03%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "19"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.352"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "29.08%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "19"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.647"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.30%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "19"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1273"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.40%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "19"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1970"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.58%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "19"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09594"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.05%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "19"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04726"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "13.92%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "19"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1048"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.16%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "19"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001324"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.28%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "19"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04198"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.22%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "19"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005905"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.73%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "19"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.348"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.19%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "19"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.405"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.17%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "19"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3507"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "4.60%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "19"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.168"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.37%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "19"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1382"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.84%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "19"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3090"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.75%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "19"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001292"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.84%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "19"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004284"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.19%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "19"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001349"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.07%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "19"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003535"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "19"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "19"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "19"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "19"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "19"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "19"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "19"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "19"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "19"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "19"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "19"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "19"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02735"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "8.87%"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "19"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05948"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "12.02%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "19"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "96.76%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "19"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.008028"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.93%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "19"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1466"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.02%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "19"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007603"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.41%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "19"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007871"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.60%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "19"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3233"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.05%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "19"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006960"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.16%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "19"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.15%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "19"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05564"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "28.07%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "19"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.14%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "19"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.15%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "19"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.15%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "19"
